$wb = $excel.ActiveWorkbook
$wsCms = $wb.Worksheets.Item("CMS")
$wsAd  = $wb.Worksheets.Item("ADCore")

# ---------------------------------------------------------------------------
# 1. CMS sheet ("Table2") gains a 5th column, "Carousel".
# ---------------------------------------------------------------------------
$tbl = $wsCms.ListObjects.Item(1)
$newCol = $tbl.ListColumns.Add()
# Setting the header cell drives both the worksheet cell and the table's
# column name/schema.
$wsCms.Range("E1").Value = "Carousel"

# Rows whose "Description" (column C) text changes from
# "Section N Accordion List" to "Section N Accordion/Custom List", and that
# also gain a value in the new "Carousel" column (E) equal to the existing
# "WXX_XX_XXX_XXX0NCL" / "WXX_XX_XXX_XXX1NCL" style codes.
$targetRows = @(24,33,42,51,60,69,78,87,96,105,114,123,132,141,150,159)

$n = 1
foreach ($r in $targetRows) {
    $wsCms.Cells.Item($r, 3).Value = "Section $n Accordion/Custom List"
    $code = "{0:D2}" -f $n
    $wsCms.Cells.Item($r, 5).Value = "WXX_XX_XXX_XXX" + $code + "CL"
    $n++
}

# Column widths: widen column C (Description) and give the new column E a
# width, mirroring the author's manual column resize / autofit.
$wsCms.Columns.Item(3).ColumnWidth = 26.42
$wsCms.Columns.Item(5).ColumnWidth = 18.42

# ---------------------------------------------------------------------------
# 2. ADCore sheet: clear the Start Date / End Date values (columns E & F)
#    for every data row, while keeping their date number format.
# ---------------------------------------------------------------------------
$wsAd.Range("E2:F186").ClearContents()

# ---------------------------------------------------------------------------
# 3. View/selection state: CMS becomes the active/selected tab with the
#    selection at F110; ADCore is no longer the selected tab, its frozen
#    top-left scroll position is reset, and its selection moves to H22.
# ---------------------------------------------------------------------------
$wsAd.Range("H22").Select()
$wsCms.Activate()
$wsCms.Range("F110").Select()
